# The commit shuffled the order of the data rows 110-119 (1-indexed sheet
# rows) in the training data sample - each row's full record (columns A:AU)
# moved to a different row position. Capture the "before" snapshot of every
# affected row first, then write each snapshot into its new destination row.
# Using a full snapshot (rather than sequential swaps) is required because
# the permutation is made of several independent cycles, not simple pairwise
# swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the full A:AU record for every row touched by the reshuffle.
$snapshot = @{}
foreach ($r in 110..119) {
    $snapshot[$r] = $ws.Range("A$r`:AU$r").Value()
}

# new row -> old row the data should come from
$mapping = @{
    110 = 114
    111 = 110
    112 = 115
    113 = 111
    114 = 113
    115 = 116
    116 = 112
    117 = 118
    118 = 119
    119 = 117
}

foreach ($r in 110..119) {
    $src = $mapping[$r]
    $ws.Range("A$r`:AU$r").Value = $snapshot[$src]
}
